$d = $word.ActiveDocument

# 1) Remove the old "_GoBack" bookmark. It currently sits, empty, in the
#    paragraph right after "Active learning ... most uncertain".
$d.Bookmarks("_GoBack").Delete()

# 2) Fix the "9constraint" typo to "(constraint" inside the
#    "Local search ..." paragraph.
$findRange = $d.Content
$findRange.Find.Execute("Local search 9constraint", $false, $false, $false, $false, $false, $true, 1, $false, "Local search (constraint", 2) | Out-Null

# 3) Locate that paragraph (by scanning Document.Paragraphs — Range.Paragraphs
#    is not reliably scoped in this host, so index from the top-level
#    collection instead).
$count = $d.Paragraphs.Count
$targetIdx = -1
for ($i = 1; $i -le $count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "Local search (constraint*") {
        $targetIdx = $i
        break
    }
}

$para = $d.Paragraphs.Item($targetIdx)
$paraText = $para.Range.Text
$paraStart = $para.Range.Start

# Split the single run into three runs at the right spots, matching the
# target structure:
#   run1: "Local search ("
#   run2: "constraint to do 5/10 grouping) - suppose a 7/8 split, "
#   <bookmarkStart _GoBack/>
#   run3: "iterating through all 8 to see if adding one to the 7 improves
#          the error score. Find clustering that gives you cleanest split"
#   <bookmarkEnd _GoBack/>   (placed right after the paragraph)
$idxConstraint = $paraText.IndexOf("constraint")
$idxIterating = $paraText.IndexOf("iterating")

$splitPos1 = $paraStart + $idxConstraint
$splitPos2 = $paraStart + $idxIterating

# Split #1: between "Local search (" and "constraint ..."
$zr1 = $d.Range($splitPos1, $splitPos1)
$d.Bookmarks.Add("ZZTMPSPLIT1", $zr1) | Out-Null
$d.Bookmarks("ZZTMPSPLIT1").Delete()

# Split #2: between "... 7/8 split, " and "iterating ..."
$zr2 = $d.Range($splitPos2, $splitPos2)
$d.Bookmarks.Add("ZZTMPSPLIT2", $zr2) | Out-Null
$d.Bookmarks("ZZTMPSPLIT2").Delete()

# 4) Re-add "_GoBack" spanning from the "iterating" split point through the
#    end of the paragraph (i.e. into the very start of the following,
#    empty paragraph) so bookmarkEnd lands right after </w:p>, as a sibling.
$nextPara = $d.Paragraphs.Item($targetIdx + 1)
$nextParaStart = $nextPara.Range.Start

$goBackRange = $d.Range($splitPos2, $nextParaStart)
$d.Bookmarks.Add("_GoBack", $goBackRange) | Out-Null

Write-Output "done"
